$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 17857692
$ws.Range("I92").Value = 18519018
$ws.Range("K92").Value = 18519018
$ws.Range("M92").Value = -18517770
$ws.Range("H99").Value = 1043.2632
$ws.Range("I99").Value = 702.13336
$ws.Range("K99").Value = 2106.40008
$ws.Range("M99").Value = -608.4000800000003
$ws.Range("H101").Value = 588.3333
$ws.Range("I101").Value = 736.25
$ws.Range("K101").Value = 2208.75
$ws.Range("M101").Value = -586.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1864.7
$ws.Range("I110").Value = 1891.36
$ws.Range("K110").Value = 1891.36
$ws.Range("M110").Value = 153.6400000000001
$ws.Range("H132").Value = 2564.1128
$ws.Range("I132").Value = 2197.2126
$ws.Range("J132").Value = 3713.7334
$ws.Range("K132").Value = 6591.6378
$ws.Range("L132").Value = 11141.2002
$ws.Range("M132").Value = -4061.6378
$ws.Range("N132").Value = -16201.2002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3497.52
$ws.Range("I105").Value = 2048.7778
$ws.Range("K105").Value = 2048.7778
$ws.Range("M105").Value = -301.7777999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 178.33333
$ws.Range("I22").Value = 178.33333
$ws.Range("K22").Value = 178.33333
$ws.Range("M22").Value = 171.66667
$ws.Range("H31").Value = 3729.0344
$ws.Range("I31").Value = 2394.0454
$ws.Range("K31").Value = 2394.0454
$ws.Range("M31").Value = -2099.0454
$ws.Range("H34").Value = 3729.0344
$ws.Range("I34").Value = 2394.0454
$ws.Range("K34").Value = 2394.0454
$ws.Range("M34").Value = -2192.0454
$ws.Range("H58").Value = 2961.9666
$ws.Range("I58").Value = 2582
$ws.Range("J58").Value = 3396.2144
$ws.Range("K58").Value = 2582
$ws.Range("L58").Value = 3396.2144
$ws.Range("M58").Value = -2379
$ws.Range("N58").Value = -3802.2144
$ws.Range("H86").Value = 3116.4
$ws.Range("I86").Value = 3089.7273
$ws.Range("J86").Value = 3189.75
$ws.Range("K86").Value = 3089.7273
$ws.Range("L86").Value = 3189.75
$ws.Range("M86").Value = -1966.7273
$ws.Range("N86").Value = -5435.75
$ws.Range("H89").Value = 3116.4
$ws.Range("I89").Value = 3089.7273
$ws.Range("J89").Value = 3189.75
$ws.Range("K89").Value = 15448.6365
$ws.Range("L89").Value = 15948.75
$ws.Range("M89").Value = -9832.636500000001
$ws.Range("N89").Value = -27180.75
$ws.Range("H107").Value = 60136.47
$ws.Range("I107").Value = 84368.914
$ws.Range("K107").Value = 84368.914
$ws.Range("M107").Value = -82448.914
$ws.Range("H132").Value = 1639.7142
$ws.Range("I132").Value = 1316.95
$ws.Range("K132").Value = 3950.85
$ws.Range("M132").Value = -1420.85
$ws.Range("H136").Value = 2961.9666
$ws.Range("I136").Value = 2582
$ws.Range("J136").Value = 3396.2144
$ws.Range("K136").Value = 7746
$ws.Range("L136").Value = 10188.6432
$ws.Range("M136").Value = -5196
$ws.Range("N136").Value = -15288.6432
$ws.Range("H141").Value = 505260.75
$ws.Range("J141").Value = 530652.0600000001
$ws.Range("L141").Value = 530652.0600000001
$ws.Range("N141").Value = -541012.0600000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 154580130
$ws.Range("I4").Value = 115139310
$ws.Range("J4").Value = 266329090
$ws.Range("K4").Value = 345417930
$ws.Range("L4").Value = 798987270
$ws.Range("M4").Value = -345417818
$ws.Range("N4").Value = -798987494
$ws.Range("H86").Value = 242
$ws.Range("I86").Value = 288
$ws.Range("K86").Value = 864
$ws.Range("M86").Value = 322
$ws.Range("H89").Value = 242
$ws.Range("I89").Value = 288
$ws.Range("K89").Value = 2592
$ws.Range("M89").Value = 3336
$ws.Range("H122").Value = 1982.5
$ws.Range("I122").Value = 1968.8334
$ws.Range("J122").Value = 2003
$ws.Range("K122").Value = 17719.5006
$ws.Range("L122").Value = 18027
$ws.Range("M122").Value = -15269.5006
$ws.Range("N122").Value = -22927
$ws.Range("H127").Value = 4273.2
$ws.Range("J127").Value = 4841.5
$ws.Range("L127").Value = 14524.5
$ws.Range("N127").Value = -24444.5
$ws.Range("H132").Value = 1779.6666
$ws.Range("J132").Value = 2670.7144
$ws.Range("L132").Value = 24036.4296
$ws.Range("N132").Value = -29096.4296

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3888.4285
$ws.Range("I122").Value = 2906.3333
$ws.Range("K122").Value = 8718.999899999999
$ws.Range("M122").Value = -6268.999899999999
$ws.Range("H132").Value = 2456.577
$ws.Range("I132").Value = 2252.75
$ws.Range("K132").Value = 6758.25
$ws.Range("M132").Value = -4228.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 7912.5
$ws.Range("J23").Value = 7900
$ws.Range("L23").Value = 7900
$ws.Range("N23").Value = -8360
$ws.Range("H43").Value = 142866860
$ws.Range("J43").Value = 166676670
$ws.Range("L43").Value = 166676670
$ws.Range("N43").Value = -166677056
$ws.Range("H139").Value = 150000
$ws.Range("J139").Value = 150000
$ws.Range("L139").Value = 150000
$ws.Range("N139").Value = -160280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 444.26315
$ws.Range("I107").Value = 280.44446
$ws.Range("K107").Value = 841.33338
$ws.Range("M107").Value = 1078.66662
$ws.Range("H115").Value = 104999
$ws.Range("J115").Value = 104999
$ws.Range("L115").Value = 104999
$ws.Range("N115").Value = -108133
$ws.Range("H122").Value = 1995
$ws.Range("J122").Value = 1995
$ws.Range("L122").Value = 5985
$ws.Range("N122").Value = -10885
$ws.Range("H127").Value = 84934.5
$ws.Range("J127").Value = 84934.5
$ws.Range("L127").Value = 84934.5
$ws.Range("N127").Value = -94854.5
$ws.Range("H132").Value = 2668.7817
$ws.Range("I132").Value = 2530.3416
$ws.Range("K132").Value = 7591.024800000001
$ws.Range("M132").Value = -5061.024800000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 29560.166
$ws.Range("I136").Value = 946.7273
$ws.Range("K136").Value = 2840.1819
$ws.Range("M136").Value = -290.1819
$ws.Range("H139").Value = 96332.336
$ws.Range("J139").Value = 96332.336
$ws.Range("L139").Value = 96332.336
$ws.Range("N139").Value = -106612.336

# ---- Remove cell N135 on WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N135").ClearContents()
